# Weekly refresh of the "Coco" price series for Vega Modelo de Temuco.
# The data table (rows 10-49) gets 5 new weekly observations woven in while
# keeping every other column (Mercado, Region, Producto, etc.) identical -
# those columns are constant for the whole table, so we only need to:
#   1) grow the used range by 5 rows (appended at the bottom, Excel copies
#      the row-above formatting down automatically - same as the diff shows
#      for the new rows' date-formatted column D),
#   2) rewrite the D/M/N/O/P/S columns for rows 10-54 with their final
#      values, and
#   3) fill in the constant columns for the five brand-new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Extend the table with 5 new rows at the bottom (A50:T54), pushing the
#    dimension from A1:T49 to A1:T54 and inheriting formatting from row 49.
$ws.Range("A50:T54").Insert()

# 2) Final per-row values for the Fecha / Volumen / Precio min / Precio max /
#    Precio promedio / Precio $/Kg columns, rows 10 through 54.
$rows = @(
    @{Row=10; D=44532; M=20;  N=28000; O=28000; P=28000; S=1400},
    @{Row=11; D=44466; M=70;  N=24000; O=24000; P=24000; S=1200},
    @{Row=12; D=44488; M=40;  N=20000; O=20000; P=20000; S=1000},
    @{Row=13; D=44425; M=15;  N=24000; O=24000; P=24000; S=1200},
    @{Row=14; D=44454; M=25;  N=25000; O=25000; P=25000; S=1250},
    @{Row=15; D=44382; M=15;  N=20000; O=20000; P=20000; S=1000},
    @{Row=16; D=44426; M=15;  N=24000; O=24000; P=24000; S=1200},
    @{Row=17; D=44421; M=20;  N=24000; O=24000; P=24000; S=1200},
    @{Row=18; D=44467; M=20;  N=24000; O=24000; P=24000; S=1200},
    @{Row=19; D=44235; M=15;  N=25000; O=25000; P=25000; S=1250},
    @{Row=20; D=44334; M=20;  N=25000; O=25000; P=25000; S=1250},
    @{Row=21; D=44356; M=15;  N=24000; O=24000; P=24000; S=1200},
    @{Row=22; D=44175; M=25;  N=23000; O=23000; P=23000; S=1150},
    @{Row=23; D=44222; M=15;  N=25000; O=25000; P=25000; S=1250},
    @{Row=24; D=44238; M=30;  N=25000; O=25000; P=25000; S=1250},
    @{Row=25; D=44468; M=20;  N=24000; O=24000; P=24000; S=1200},
    @{Row=26; D=44442; M=25;  N=23000; O=23000; P=23000; S=1150},
    @{Row=27; D=44214; M=15;  N=25000; O=25000; P=25000; S=1250},
    @{Row=28; D=44194; M=20;  N=20000; O=20000; P=20000; S=1000},
    @{Row=29; D=44389; M=20;  N=24000; O=24000; P=24000; S=1200},
    @{Row=30; D=44349; M=30;  N=24000; O=24000; P=24000; S=1200},
    @{Row=31; D=44412; M=20;  N=25000; O=25000; P=25000; S=1250},
    @{Row=32; D=44398; M=15;  N=25000; O=25000; P=25000; S=1250},
    @{Row=33; D=44420; M=35;  N=25000; O=25000; P=25000; S=1250},
    @{Row=34; D=44249; M=15;  N=25000; O=25000; P=25000; S=1250},
    @{Row=35; D=44232; M=15;  N=25000; O=25000; P=25000; S=1250},
    @{Row=36; D=44431; M=40;  N=24000; O=24000; P=24000; S=1200},
    @{Row=37; D=44400; M=5;   N=24000; O=24000; P=24000; S=1200},
    @{Row=38; D=44363; M=30;  N=24000; O=24000; P=24000; S=1200},
    @{Row=39; D=44390; M=10;  N=24000; O=24000; P=24000; S=1200},
    @{Row=40; D=44461; M=30;  N=24000; O=24000; P=24000; S=1200},
    @{Row=41; D=44221; M=30;  N=25000; O=25000; P=25000; S=1250},
    @{Row=42; D=44432; M=30;  N=24000; O=24000; P=24000; S=1200},
    @{Row=43; D=44428; M=15;  N=24000; O=24000; P=24000; S=1200},
    @{Row=44; D=44462; M=10;  N=24000; O=24000; P=24000; S=1200},
    @{Row=45; D=44435; M=100; N=24000; O=24000; P=24000; S=1200},
    @{Row=46; D=44231; M=15;  N=25000; O=25000; P=25000; S=1250},
    @{Row=47; D=44489; M=40;  N=24000; O=24000; P=24000; S=1200},
    @{Row=48; D=44377; M=15;  N=20000; O=20000; P=20000; S=1000},
    @{Row=49; D=44396; M=12;  N=24000; O=24000; P=24000; S=1200},
    @{Row=50; D=44418; M=20;  N=24000; O=24000; P=24000; S=1200},
    @{Row=51; D=44469; M=40;  N=24000; O=24000; P=24000; S=1200},
    @{Row=52; D=44392; M=10;  N=24000; O=24000; P=24000; S=1200},
    @{Row=53; D=44474; M=20;  N=24000; O=24000; P=24000; S=1200},
    @{Row=54; D=44424; M=25;  N=24000; O=24000; P=24000; S=1200}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = $r.D    # D - Fecha
    $ws.Cells.Item($r.Row, 13).Value = $r.M   # M - Volumen
    $ws.Cells.Item($r.Row, 14).Value = $r.N   # N - Precio minimo
    $ws.Cells.Item($r.Row, 15).Value = $r.O   # O - Precio maximo
    $ws.Cells.Item($r.Row, 16).Value = $r.P   # P - Precio promedio ponderado
    $ws.Cells.Item($r.Row, 19).Value = $r.S   # S - Precio $/Kg
}

# 3) The five appended rows (50-54) were entirely blank after the insert;
#    fill in the columns that are constant across the whole table.
for ($r = 50; $r -le 54; $r++) {
    $ws.Cells.Item($r, 1).Value = 10                               # A - Mercado ID
    $ws.Cells.Item($r, 2).Value = "Vega Modelo de Temuco"           # B - Mercado
    $ws.Cells.Item($r, 3).Value = "La Araucanía"                   # C - Region
    $ws.Cells.Item($r, 5).Value = 9                                 # E - Codreg
    $ws.Cells.Item($r, 6).Value = "Fruta"                           # F - Tipo
    $ws.Cells.Item($r, 7).Value = 100108                            # G - Producto ID
    $ws.Cells.Item($r, 8).Value = "Tropicales y subtropicales"      # H - Producto
    $ws.Cells.Item($r, 9).Value = 100108007                         # I - Categoria ID
    $ws.Cells.Item($r, 10).Value = "Coco"                           # J - Categoria
    $ws.Cells.Item($r, 11).Value = "Sin especificar"                # K - Variedad
    $ws.Cells.Item($r, 12).Value = "Primera"                        # L - Calidad
    $ws.Cells.Item($r, 17).Value = "`$/malla 20 unidades"           # Q - Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value = "Perú"                           # R - Origen
    $ws.Cells.Item($r, 20).Value = 20                               # T - Kg / unidad
}
